$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: date (06/10/2013) in A7, value 6 in B7
# Copy the date formatting from an existing date cell so the same
# (built-in) number-format style is reused rather than creating a new one.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = 41553
$ws.Range("B7").Value = 6

# Update the active selection to match the authored state
$ws.Range("C10").Select()
